$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: was a merged/centered C9:E9 block, becomes a plain 0/1/2 header row ---
$ws.Range("C9:E9").UnMerge()
$ws.Range("C9:E9").ClearFormats()
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2

# --- Row 10: row index 0 + 2D array header values (strings), styled like the
#     "top of box" cells (green fill, border all sides except bottom) plus a
#     trailing "no border" green cell ---
$ws.Range("B10").Value = 0

$ws.Range("C6").Copy()
$ws.Range("C10:E10").PasteSpecial(-4122)  # xlPasteFormats (fillId3 borderId1 full box like a closed cell)

$ws.Range("C3").Copy()
$ws.Range("C10:E10").PasteSpecial(-4122)  # fillId3 borderId2 (top/left/right, no bottom)

$ws.Range("C10").Value = "Delhi"
$ws.Range("D10").Value = "Mumbai"
$ws.Range("E10").Value = "UP"

$ws.Range("F10").Interior.Color = 5296274
$ws.Range("F10").Value = "Bihar"

# --- Row 11: row index 1 + styled cells continuing the box, with a special
#     "code font" cell in the first column ---
$ws.Range("B11").Value = 1

$ws.Range("C6").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Font.Name = "Consolas"
$ws.Range("C11").Font.Size = 8
$ws.Range("C11").Font.Color = 7901646
$ws.Range("C11").VerticalAlignment = -4108  # xlCenter

$ws.Range("C3").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("E11:F11").PasteSpecial(-4122)

# --- Row 12: row index 2 + bottom of the box ---
$ws.Range("B12").Value = 2

$ws.Range("C3").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)

$ws.Range("E12").Interior.Color = 5296274

$excel.CutCopyMode = 0

# --- View: scroll down a bit and move the selection ---
$ws.Range("G9").Select()

# --- Page setup: portrait orientation (as in target) ---
$ws.PageSetup.Orientation = 1
